$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the trailing bold "Play Atlantis Megaways Free | Online
#    Slot Game Review" paragraph (it now lives near the top instead,
#    as plain Heading1 text -- this one is the bold *body* copy).
# ------------------------------------------------------------------
$titleText = "Play Atlantis Megaways Free | Online Slot Game Review"
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    $t = $t.TrimEnd([char]13, [char]7)
    if ($t -eq $titleText -and $p.Style.NameLocal -ne "Heading 1") {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2) Replace the final paragraph's text (the old meta-description
#    sentence) with the new "Create a cartoon-style feature image..."
#    image-prompt text, keeping its italic run formatting intact.
# ------------------------------------------------------------------
$oldTail = "Read our review of Atlantis Megaways, an underwater-themed online slot with cascading reels, 4 jackpot prizes, and up to 117,649 ways to win. Play for free today."
$newTail = "Create a cartoon-style feature image for Atlantis Megaways that features a happy Maya warrior wearing glasses. The warrior should be positioned underwater among ruins of the lost city of Atlantis with sea creatures swimming around in the background. The image should incorporate the game's logo and feature vibrant colors that capture the adventurous and mysterious theme of the game. The image should also clearly convey the idea of winning cash prizes with a bubbly, celebratory vibe."

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$pr = $lastPara.Range
$textRange = $d.Range($pr.Start, $pr.End - 1)
if ($textRange.Text -eq $oldTail) {
    $textRange.Text = $newTail
}

# ------------------------------------------------------------------
# 3) Insert a new paragraph right after the (now only) top heading
#    paragraph: "Meta description" (bold) followed by the plain-text
#    meta description sentence.
# ------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$boldLabel = "Meta description"
$restText = ": Read our review of Atlantis Megaways, an underwater-themed online slot with cascading reels, 4 jackpot prizes, and up to 117,649 ways to win. Play for free today."

$metaRangeFull = $metaPara.Range
$metaRange = $d.Range($metaRangeFull.Start, $metaRangeFull.End - 1)
$metaRange.Text = $boldLabel + $restText

$boldRange = $d.Range($metaRange.Start, $metaRange.Start + $boldLabel.Length)
$boldRange.Font.Bold = 1

Write-Output "done"
